# Adds the `shellyrelay` IO Command row to the "Commands" sheet, matching
# the upstream commit: "adds `shellyrelay` IO Command to switch Shelly
# Plugs on/off - updated libs - updated translations"
#
# A new row is inserted right above the existing "S7 Command" section
# (old row 71, now row 72), containing the command signature in column B
# and its description in column C - following the same layout used by
# every other command row in this table.

$wb = $excel.ActiveWorkbook
$wsCommands = $wb.Worksheets.Item("Commands")

$wsCommands.Activate()

# Insert a new row at 71; existing rows 71.. shift down to 72..
$wsCommands.Rows.Item(71).Insert()

# Fill in the new command row (column B = command signature, column C = description)
$wsCommands.Range("B71").Value2 = "shellyrelay(n,b)"
$wsCommands.Range("C71").Value2 = "switches Shelly plug number <n> ON if b is true or 1, and OFF otherwise"

# Leave the same cells selected, as in the saved workbook
$wsCommands.Range("B71:C71").Select()
